# GMS Data Release 1
# The "observation" data-dictionary row for patient_id is renamed to
# participant_id (terminology change from "patient" to "participant").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "participant_id"

# Leave the selection where the author left it after making the edit.
$ws.Range("F5").Select()
